$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 125002504
$ws.Range("I86").Value = 500000640
$ws.Range("J86").Value = 3124.5
$ws.Range("K86").Value = 500000640
$ws.Range("L86").Value = 3124.5
$ws.Range("M86").Value = -499999517
$ws.Range("N86").Value = -5370.5
$ws.Range("H89").Value = 125002504
$ws.Range("I89").Value = 500000640
$ws.Range("J89").Value = 3124.5
$ws.Range("K89").Value = 2500003200
$ws.Range("L89").Value = 15622.5
$ws.Range("M89").Value = -2499997584
$ws.Range("N89").Value = -26854.5
$ws.Range("H100").Value = 3247.25
$ws.Range("I100").Value = 2996.3333
$ws.Range("K100").Value = 2996.3333
$ws.Range("M100").Value = -2455.3333
$ws.Range("H106").Value = 3283
$ws.Range("I106").Value = 3566
$ws.Range("J106").Value = 3000
$ws.Range("K106").Value = 3566
$ws.Range("L106").Value = 3000
$ws.Range("M106").Value = -2935
$ws.Range("N106").Value = -4262
$ws.Range("H113").Value = 5100.5
$ws.Range("J113").Value = 5433.3335
$ws.Range("L113").Value = 5433.3335
$ws.Range("N113").Value = -11941.3335
$ws.Range("H133").Value = 105996
$ws.Range("J133").Value = 105996
$ws.Range("L133").Value = 105996
$ws.Range("N133").Value = -116116
$ws.Range("H137").Value = 3721.6743
$ws.Range("I137").Value = 1591.3529
$ws.Range("K137").Value = 4774.0587
$ws.Range("M137").Value = -2224.0587
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 46668.363
$ws.Range("I45").Value = 67910.42999999999
$ws.Range("J45").Value = 9494.75
$ws.Range("K45").Value = 67910.42999999999
$ws.Range("L45").Value = 9494.75
$ws.Range("M45").Value = -67533.42999999999
$ws.Range("N45").Value = -10248.75
$ws.Range("H61").Value = 3964.739
$ws.Range("I61").Value = 2184.9333
$ws.Range("J61").Value = 7301.875
$ws.Range("K61").Value = 2184.9333
$ws.Range("L61").Value = 7301.875
$ws.Range("M61").Value = -1972.9333
$ws.Range("N61").Value = -7725.875
$ws.Range("H132").Value = 2193.5898
$ws.Range("I132").Value = 1628.4073
$ws.Range("K132").Value = 4885.2219
$ws.Range("M132").Value = -2355.2219
$ws.Range("H136").Value = 3964.739
$ws.Range("I136").Value = 2184.9333
$ws.Range("J136").Value = 7301.875
$ws.Range("K136").Value = 6554.7999
$ws.Range("L136").Value = 21905.625
$ws.Range("M136").Value = -4004.7999
$ws.Range("N136").Value = -27005.625
$ws.Range("H11").Value = 289.5
$ws.Range("I11").Value = 504
$ws.Range("J11").Value = 75
$ws.Range("K11").Value = 504
$ws.Range("L11").Value = 75
$ws.Range("M11").Value = -364
$ws.Range("N11").Value = -355
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 204501.8
$ws.Range("I99").Value = 501754.5
$ws.Range("K99").Value = 501754.5
$ws.Range("M99").Value = -500256.5
$ws.Range("H134").Value = 2744.818
$ws.Range("I134").Value = 1962.2142
$ws.Range("K134").Value = 5886.642599999999
$ws.Range("M134").Value = -3351.642599999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").ClearContents()
$ws.Range("N12").Value = 0
$ws.Range("H22").Value = 1824.2727
$ws.Range("I22").Value = 1739.2
$ws.Range("J22").Value = 1895.1666
$ws.Range("K22").Value = 1739.2
$ws.Range("L22").Value = 1895.1666
$ws.Range("M22").Value = -1389.2
$ws.Range("N22").Value = -2595.1666
$ws.Range("H31").Value = 3848.0312
$ws.Range("I31").Value = 3161.1035
$ws.Range("J31").Value = 10488.333
$ws.Range("K31").Value = 3161.1035
$ws.Range("L31").Value = 10488.333
$ws.Range("M31").Value = -2866.1035
$ws.Range("N31").Value = -11078.333
$ws.Range("H34").Value = 3848.0312
$ws.Range("I34").Value = 3161.1035
$ws.Range("J34").Value = 10488.333
$ws.Range("K34").Value = 3161.1035
$ws.Range("L34").Value = 10488.333
$ws.Range("M34").Value = -2959.1035
$ws.Range("N34").Value = -10892.333
$ws.Range("H58").Value = 3800.3684
$ws.Range("I58").Value = 2994.9
$ws.Range("J58").Value = 4695.3335
$ws.Range("K58").Value = 2994.9
$ws.Range("L58").Value = 4695.3335
$ws.Range("M58").Value = -2791.9
$ws.Range("N58").Value = -5101.3335
$ws.Range("H105").Value = 3209
$ws.Range("J105").Value = 3209
$ws.Range("L105").Value = 3209
$ws.Range("N105").Value = -6703
$ws.Range("H132").Value = 12197772
$ws.Range("I132").Value = 13515814
$ws.Range("J132").Value = 5874.75
$ws.Range("K132").Value = 40547442
$ws.Range("L132").Value = 17624.25
$ws.Range("M132").Value = -40544912
$ws.Range("N132").Value = -22684.25
$ws.Range("H134").Value = 3023.9412
$ws.Range("I134").Value = 2775.4375
$ws.Range("J134").Value = 7000
$ws.Range("K134").Value = 8326.3125
$ws.Range("L134").Value = 21000
$ws.Range("M134").Value = -5791.3125
$ws.Range("N134").Value = -26070
$ws.Range("H136").Value = 3800.3684
$ws.Range("I136").Value = 2994.9
$ws.Range("J136").Value = 4695.3335
$ws.Range("K136").Value = 8984.700000000001
$ws.Range("L136").Value = 14086.0005
$ws.Range("M136").Value = -6434.700000000001
$ws.Range("N136").Value = -19186.0005
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3087.7827
$ws.Range("J5").Value = 4167.6875
$ws.Range("L5").Value = 12503.0625
$ws.Range("N5").Value = -12727.0625
$ws.Range("H22").Value = 76928220
$ws.Range("J22").Value = 6586
$ws.Range("L22").Value = 19758
$ws.Range("N22").Value = -20096
$ws.Range("H27").Value = 76928220
$ws.Range("J27").Value = 6586
$ws.Range("L27").Value = 19758
$ws.Range("N27").Value = -19962
$ws.Range("H129").Value = 2180.8572
$ws.Range("I129").Value = 2206.4285
$ws.Range("K129").Value = 6619.2855
$ws.Range("M129").Value = -1619.2855
$ws.Range("H135").Value = 3087.7827
$ws.Range("J135").Value = 4167.6875
$ws.Range("L135").Value = 37509.1875
$ws.Range("N135").Value = -42579.1875
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 58825030
$ws.Range("I80").Value = 125001384
$ws.Range("J80").Value = 1605.4445
$ws.Range("K80").Value = 125001384
$ws.Range("L80").Value = 1605.4445
$ws.Range("M80").Value = -125000386
$ws.Range("N80").Value = -3601.4445
$ws.Range("H83").Value = 58825030
$ws.Range("I83").Value = 125001384
$ws.Range("J83").Value = 1605.4445
$ws.Range("K83").Value = 625006920
$ws.Range("L83").Value = 8027.2225
$ws.Range("M83").Value = -625001928
$ws.Range("N83").Value = -18011.2225
$ws.Range("H102").Value = 4732.3774
$ws.Range("I102").Value = 788.6591
$ws.Range("J102").Value = 24012.777
$ws.Range("K102").Value = 788.6591
$ws.Range("L102").Value = 24012.777
$ws.Range("M102").Value = 833.3409
$ws.Range("N102").Value = -27256.777
$ws.Range("H132").Value = 4521.7856
$ws.Range("I132").Value = 4730.5
$ws.Range("K132").Value = 14191.5
$ws.Range("M132").Value = -11661.5
$ws.Range("H133").Value = 165000
$ws.Range("J133").Value = 165000
$ws.Range("L133").Value = 165000
$ws.Range("N133").Value = -175120
$ws.Range("H135").Value = 68574.63
$ws.Range("J135").Value = 68574.63
$ws.Range("L135").Value = 68574.63
$ws.Range("N135").Value = -78714.63
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 844.6429000000001
$ws.Range("I16").Value = 1034.4
$ws.Range("K16").Value = 1034.4
$ws.Range("M16").Value = -864.4000000000001
$ws.Range("H61").Value = 1141.2222
$ws.Range("I61").Value = 969.8333
$ws.Range("K61").Value = 969.8333
$ws.Range("M61").Value = -767.8333
$ws.Range("H68").Value = 7501.1665
$ws.Range("I68").Value = 7000.5
$ws.Range("J68").Value = 7751.5
$ws.Range("K68").Value = 7000.5
$ws.Range("L68").Value = 7751.5
$ws.Range("M68").Value = -6251.5
$ws.Range("N68").Value = -9249.5
$ws.Range("H71").Value = 7501.1665
$ws.Range("I71").Value = 7000.5
$ws.Range("J71").Value = 7751.5
$ws.Range("K71").Value = 35002.5
$ws.Range("L71").Value = 38757.5
$ws.Range("M71").Value = -31258.5
$ws.Range("N71").Value = -46245.5
$ws.Range("H113").Value = 1141.2222
$ws.Range("I113").Value = 969.8333
$ws.Range("K113").Value = 969.8333
$ws.Range("M113").Value = 1200.1667
$ws.Range("H122").Value = 2986.1538
$ws.Range("J122").Value = 2514.4285
$ws.Range("L122").Value = 7543.2855
$ws.Range("N122").Value = -12443.2855
$ws.Range("H136").Value = 4118.484
$ws.Range("I136").Value = 3845.9048
$ws.Range("K136").Value = 11537.7144
$ws.Range("M136").Value = -8987.714399999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 25000
$ws.Range("J52").Value = 25000
$ws.Range("L52").Value = 25000
$ws.Range("N52").Value = -25452
